$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-12-27 Wednesday" "2023-12-28 Thursday"

Replace-Text "50÷9=5, 5" "42÷6=7, 0"
Replace-Text "24÷7=3, 3" "47÷5=9, 2"
Replace-Text "75÷2=37, 1" "79÷5=15, 4"
Replace-Text "64÷6=10, 4" "46÷6=7, 4"
Replace-Text "46÷8=5, 6" "37÷2=18, 1"

Replace-Text "47÷3=15, 2" "14÷9=1, 5"
Replace-Text "43÷9=4, 7" "44÷9=4, 8"
Replace-Text "50÷3=16, 2" "44÷4=11, 0"
Replace-Text "68÷8=8, 4" "42÷3=14, 0"
Replace-Text "94÷8=11, 6" "15÷6=2, 3"

Replace-Text "72÷6=12, 0" "86÷4=21, 2"
Replace-Text "96÷7=13, 5" "33÷3=11, 0"
Replace-Text "17÷5=3, 2" "56÷8=7, 0"
Replace-Text "40÷5=8, 0" "83÷2=41, 1"
Replace-Text "25÷7=3, 4" "93÷7=13, 2"

Replace-Text "19÷3=6, 1" "42÷9=4, 6"
Replace-Text "11÷6=1, 5" "12÷9=1, 3"
Replace-Text "40÷3=13, 1" "64÷6=10, 4"
Replace-Text "74÷4=18, 2" "56÷9=6, 2"
Replace-Text "27÷4=6, 3" "42÷5=8, 2"

Replace-Text "28÷8=3, 4" "18÷4=4, 2"
Replace-Text "24÷3=8, 0" "20÷7=2, 6"
Replace-Text "98÷5=19, 3" "75÷9=8, 3"
Replace-Text "26÷4=6, 2" "17÷5=3, 2"
Replace-Text "25÷6=4, 1" "37÷6=6, 1"

Write-Output "Done"
